$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D=44260; J=400; K=37000; L=38000; M=37500; P=2885},
    @{Row=3; D=44312; J=400; K=26000; L=27000; M=26500; P=2038},
    @{Row=4; D=44428; J=480; K=14000; L=15000; M=14500; P=1115},
    @{Row=5; D=44419; J=600; K=14000; L=15000; M=14500; P=1115},
    @{Row=6; D=44498; J=400; K=14000; L=15000; M=14500; P=1115},
    @{Row=7; D=44410; J=600; K=14000; L=15000; M=14500; P=1115},
    @{Row=8; D=44414; J=500; K=14000; L=15000; M=14500; P=1115},
    @{Row=9; D=44249; J=400; K=42000; L=43000; M=42500; P=3269},
    @{Row=10; D=44365; J=500; K=19500; L=20000; M=19750; P=1519},
    @{Row=11; D=44426; J=460; K=14000; L=15000; M=14500; P=1115},
    @{Row=12; D=44344; J=400; K=18500; L=19000; M=18750; P=1442},
    @{Row=13; D=44435; J=480; K=13000; L=14000; M=13500; P=1038},
    @{Row=14; D=44379; J=600; K=17000; L=18000; M=17500; P=1346},
    @{Row=15; D=44309; J=400; K=26000; L=27000; M=26500; P=2038},
    @{Row=16; D=44412; J=600; K=14000; L=15000; M=14500; P=1115},
    @{Row=17; D=44335; J=480; K=24500; L=25000; M=24750; P=1904},
    @{Row=18; D=44326; J=460; K=25000; L=26000; M=25500; P=1962},
    @{Row=19; D=44418; J=500; K=14000; L=15000; M=14500; P=1115},
    @{Row=20; D=44242; J=400; K=44000; L=45000; M=44500; P=3423},
    @{Row=21; D=44442; J=460; K=14000; L=15000; M=14500; P=1115},
    @{Row=22; D=44333; J=440; K=24000; L=25000; M=24500; P=1885},
    @{Row=23; D=44323; J=460; K=25000; L=26000; M=25500; P=1962},
    @{Row=24; D=44400; J=600; K=15000; L=16000; M=15500; P=1192},
    @{Row=25; D=44484; J=360; K=14000; L=15000; M=14500; P=1115},
    @{Row=26; D=44445; J=600; K=13000; L=14000; M=13500; P=1038},
    @{Row=27; D=44383; J=200; K=17000; L=18000; M=17500; P=1346},
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 16).Value = $row.P
}
